$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers
$ws.Range("A1").Value = "sno"
$ws.Range("B1").Value = "fname"
$ws.Range("C1").Value = "lname"
$ws.Range("D1").Value = "company"
$ws.Range("E1").Value = "address"
$ws.Range("F1").Value = "address2"

# Row 2: data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "maha"
$ws.Range("C2").Value = "m"
$ws.Range("D2").Value = "mahagroup"
$ws.Range("E2").Value = "toraipakkam"
$ws.Range("F2").Value = "greenstechnology"

$ws.Range("F2").Select() | Out-Null
